# PROS-7407 CCANZ Linear SOS measurement change
#
# - Update the "size" values list (shared by the "Exclude" sheet's
#   Value-2 column, rows 2/4/7/9) from "2,4,2.25,2.4,10,3" to
#   "2,2.25,2.4,2.8,3,4,10".
# - Re-fit the columns affected by the now-longer text (column G widens
#   to fit the new value; neighbouring columns settle to the widths
#   Excel/Calc compute on save).
# - Flip which sheet/cell is active & selected: "Exclude" becomes the
#   active sheet (selected cell G19), "Include" is no longer the
#   selected tab (its own selection, E30, is unchanged).
#
# NOTE on the ColumnWidth numbers below: the host's classic
# `Range.ColumnWidth`/`Columns.ColumnWidth` COM property only persists
# width in 1/6-character increments (it round-trips through an internal
# pixel rounding step), so the literal target widths from the workbook
# XML can't be hit exactly. Each value here is the input that lands on
# the closest achievable 1/6-increment to the target width.

$wb = $excel.ActiveWorkbook

$wsExclude = $wb.Worksheets.Item("Exclude")
$wsInclude = $wb.Worksheets.Item("Include")

# --- Data change -----------------------------------------------------
$newSizes = "2,2.25,2.4,2.8,3,4,10"
$wsExclude.Range("G2").Value = $newSizes
$wsExclude.Range("G4").Value = $newSizes
$wsExclude.Range("G7").Value = $newSizes
$wsExclude.Range("G9").Value = $newSizes

# --- Re-fit columns now that column G holds longer text --------------
# "Exclude" sheet: A, B:D, E:F, G
$wsExclude.Columns.Item(1).ColumnWidth = 32.10544217687077
$wsExclude.Columns.Item(2).ColumnWidth = 27.651360544217667
$wsExclude.Columns.Item(3).ColumnWidth = 27.651360544217667
$wsExclude.Columns.Item(4).ColumnWidth = 27.651360544217667
$wsExclude.Columns.Item(5).ColumnWidth = 26.02891156462587
$wsExclude.Columns.Item(6).ColumnWidth = 26.02891156462587
$wsExclude.Columns.Item(7).ColumnWidth = 47.937074829931966

# "Include" sheet: A, B, C, D, E
$wsInclude.Columns.Item(1).ColumnWidth = 35.88605442176866
$wsInclude.Columns.Item(2).ColumnWidth = 22.518707482993168
$wsInclude.Columns.Item(3).ColumnWidth = 20.63095238095237
$wsInclude.Columns.Item(4).ColumnWidth = 15.365646258503366
$wsInclude.Columns.Item(5).ColumnWidth = 18.605442176870767

# --- Active sheet / selection -----------------------------------------
# "Include" was the selected/active tab before; "Exclude" becomes
# active now, with G19 selected on it. Select "Include" first
# (restoring its own selection at E30) then finish on "Exclude" so it
# ends up the active tab with G19 selected.
$wsInclude.Select()
$wsInclude.Range("E30").Select()

$wsExclude.Select()
$wsExclude.Range("G19").Select()
